$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrement the index column (B1:B44) by 1 -> new sequence 0..43
for ($r = 1; $r -le 44; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 1
}

# Clear the stray phone-number value in C9 (keeps its existing style/formatting)
$ws.Range("C9").ClearContents()

# Update the selected range to just C9
$ws.Range("C9").Select()
